$d = $word.ActiveDocument

# The page (docs/EF/lom3244.docx) was rebuilt by the static-site generator
# and no longer scrapes the "Ver no Jupiter / Salvar em pdf / Salvar em
# docx" toolbar line nor the "(c) 2020 ... Jekyll and Github pages ..."
# footer line that used to trail the "Requisitos" section. Remove the
# blank spacer paragraph plus those two text paragraphs, which sit right
# after the "LOM3206: Eletronica (Requisito)" line and right before the
# paragraph carrying the trailing page break.

$anchorIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*LOM3206*Eletr*nica*Requisito*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -gt 0) {
    # Deleting index ($anchorIndex + 1) three times removes the blank
    # paragraph, the "Ver no Jupiter ..." paragraph and the "(c) 2020 ..."
    # paragraph in turn, since each deletion shifts the following
    # paragraphs up into that same slot.
    $target = $anchorIndex + 1
    $d.Paragraphs.Item($target).Range.Delete()
    $d.Paragraphs.Item($target).Range.Delete()
    $d.Paragraphs.Item($target).Range.Delete()
}

Write-Output ("anchorIndex=" + $anchorIndex)
Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
